$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A180").Value = "What are the size of curve units?"
$ws.Range("B180").Value = "llama3.2:latest"
$ws.Range("C180").Value = "The size of curve units is 24."

$ws.Range("A181").Value = "How many queries can be defined per zone type?"
$ws.Range("B181").Value = "llama3.2:latest"
$ws.Range("C181").Value = "According to the document, there are 75 query definitions per zone type."
